$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = 44260
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 70
$ws.Range("K2").Value = 16000
$ws.Range("L2").Value = 16000
$ws.Range("M2").Value = 16000
$ws.Range("P2").Value = 889

$ws.Range("D3").Value2 = 44312
$ws.Range("J3").Value = 15

$ws.Range("D4").Value2 = 44250
$ws.Range("J4").Value = 60
$ws.Range("K4").Value = 18000
$ws.Range("L4").Value = 18000
$ws.Range("M4").Value = 18000
$ws.Range("P4").Value = 1000

$ws.Range("D5").Value2 = 44236
$ws.Range("I5").Value = "Especial"
$ws.Range("J5").Value = 60
$ws.Range("K5").Value = 20000
$ws.Range("L5").Value = 20000
$ws.Range("M5").Value = 20000
$ws.Range("P5").Value = 1111

$ws.Range("D6").Value2 = 44253
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 70

$ws.Range("D7").Value2 = 44232
$ws.Range("I7").Value = "Especial"
$ws.Range("J7").Value = 50
$ws.Range("K7").Value = 22000
$ws.Range("L7").Value = 22000
$ws.Range("M7").Value = 22000
$ws.Range("P7").Value = 1222

$ws.Range("D8").Value2 = 44239
$ws.Range("J8").Value = 60
$ws.Range("K8").Value = 20000
$ws.Range("L8").Value = 20000
$ws.Range("M8").Value = 20000
$ws.Range("P8").Value = 1111

$ws.Range("D9").Value2 = 44285
$ws.Range("J9").Value = 70

$ws.Range("D10").Value2 = 44243
$ws.Range("I10").Value = "Especial"
$ws.Range("J10").Value = 60
$ws.Range("K10").Value = 20000
$ws.Range("L10").Value = 20000
$ws.Range("M10").Value = 20000
$ws.Range("P10").Value = 1111

$ws.Range("D11").Value2 = 44267
$ws.Range("I11").Value = "Primera"
$ws.Range("K11").Value = 16000
$ws.Range("L11").Value = 16000
$ws.Range("M11").Value = 16000
$ws.Range("P11").Value = 889

$ws.Range("D12").Value2 = 44264
$ws.Range("I12").Value = "Primera"
$ws.Range("J12").Value = 80
$ws.Range("K12").Value = 16000
$ws.Range("L12").Value = 16000
$ws.Range("M12").Value = 16000
$ws.Range("P12").Value = 889

$ws.Range("D13").Value2 = 44271
$ws.Range("J13").Value = 70

$ws.Range("D14").Value2 = 44299
$ws.Range("J14").Value = 80
$ws.Range("K14").Value = 18000
$ws.Range("L14").Value = 18000
$ws.Range("M14").Value = 18000
$ws.Range("P14").Value = 1000

$ws.Range("D15").Value2 = 44320
$ws.Range("J15").Value = 90
$ws.Range("K15").Value = 17500
$ws.Range("L15").Value = 17500
$ws.Range("M15").Value = 17500
$ws.Range("P15").Value = 972

$ws.Range("D16").Value2 = 44364
$ws.Range("I16").Value = "Especial"
$ws.Range("J16").Value = 30
$ws.Range("K16").Value = 20000
$ws.Range("L16").Value = 20000
$ws.Range("M16").Value = 20000
$ws.Range("P16").Value = 1111

$ws.Range("D17").Value2 = 44365
$ws.Range("I17").Value = "Especial"
$ws.Range("J17").Value = 50
$ws.Range("K17").Value = 20000
$ws.Range("L17").Value = 20000
$ws.Range("M17").Value = 20000
$ws.Range("P17").Value = 1111

$ws.Range("D18").Value2 = 44257
$ws.Range("K18").Value = 16000
$ws.Range("L18").Value = 16000
$ws.Range("M18").Value = 16000
$ws.Range("P18").Value = 889

$ws.Range("D19").Value2 = 44252
$ws.Range("I19").Value = "Primera"
$ws.Range("J19").Value = 40
$ws.Range("K19").Value = 18000
$ws.Range("L19").Value = 18000
$ws.Range("M19").Value = 18000
$ws.Range("P19").Value = 1000

$ws.Range("D20").Value2 = 44259
$ws.Range("I20").Value = "Primera"
$ws.Range("J20").Value = 70

$ws.Range("D21").Value2 = 44313
$ws.Range("J21").Value = 80
$ws.Range("K21").Value = 18000
$ws.Range("L21").Value = 18000
$ws.Range("M21").Value = 18000
$ws.Range("P21").Value = 1000

$ws.Range("D22").Value2 = 44251
$ws.Range("I22").Value = "Primera"
$ws.Range("J22").Value = 20

$ws.Range("D23").Value2 = 44327
$ws.Range("I23").Value = "Especial"
$ws.Range("J23").Value = 80
$ws.Range("K23").Value = 16000
$ws.Range("L23").Value = 16000
$ws.Range("M23").Value = 16000
$ws.Range("P23").Value = 889

$ws.Range("D24").Value2 = 44245
$ws.Range("I24").Value = "Primera"
$ws.Range("J24").Value = 40

$ws.Range("D25").Value2 = 44309
$ws.Range("I25").Value = "Especial"
$ws.Range("K25").Value = 18000
$ws.Range("L25").Value = 18000
$ws.Range("M25").Value = 18000
$ws.Range("P25").Value = 1000

$ws.Range("D26").Value2 = 44238
$ws.Range("K26").Value = 20000
$ws.Range("L26").Value = 20000
$ws.Range("M26").Value = 20000
$ws.Range("P26").Value = 1111

$ws.Range("D27").Value2 = 44274
$ws.Range("I27").Value = "Primera"
$ws.Range("K27").Value = 16000
$ws.Range("L27").Value = 16000
$ws.Range("M27").Value = 16000
$ws.Range("P27").Value = 889

$ws.Range("D28").Value2 = 44326
$ws.Range("J28").Value = 15
$ws.Range("K28").Value = 18000
$ws.Range("L28").Value = 18000
$ws.Range("M28").Value = 18000
$ws.Range("P28").Value = 1000

$ws.Range("D29").Value2 = 44350
$ws.Range("I29").Value = "Primera"
$ws.Range("J29").Value = 20
$ws.Range("K29").Value = 20000
$ws.Range("L29").Value = 20000
$ws.Range("M29").Value = 20000
$ws.Range("P29").Value = 1111

$ws.Range("D30").Value2 = 44278
$ws.Range("J30").Value = 70

$ws.Range("D31").Value2 = 44371
$ws.Range("I31").Value = "Especial"
$ws.Range("J31").Value = 20
$ws.Range("K31").Value = 20000
$ws.Range("L31").Value = 20000
$ws.Range("M31").Value = 20000
$ws.Range("P31").Value = 1111

$ws.Range("D32").Value2 = 44357
$ws.Range("J32").Value = 15
$ws.Range("K32").Value = 20000
$ws.Range("L32").Value = 20000
$ws.Range("M32").Value = 20000
$ws.Range("P32").Value = 1111

$ws.Range("D33").Value2 = 44316
$ws.Range("I33").Value = "Especial"
$ws.Range("J33").Value = 70

$ws.Range("D34").Value2 = 44242
$ws.Range("I34").Value = "Especial"
$ws.Range("J34").Value = 50

$ws.Range("D35").Value2 = 44280
$ws.Range("J35").Value = 40
$ws.Range("K35").Value = 18000
$ws.Range("L35").Value = 18000
$ws.Range("M35").Value = 18000
$ws.Range("P35").Value = 1000

$ws.Range("D36").Value2 = 44292
$ws.Range("J36").Value = 70
$ws.Range("K36").Value = 17000
$ws.Range("L36").Value = 17000
$ws.Range("M36").Value = 17000
$ws.Range("P36").Value = 944

$ws.Range("D37").Value2 = 44323
$ws.Range("J37").Value = 70

$ws.Range("D38").Value2 = 44306
$ws.Range("J38").Value = 80

$ws.Range("D39").Value2 = 44246
$ws.Range("I39").Value = "Primera"
$ws.Range("J39").Value = 60
$ws.Range("K39").Value = 18000
$ws.Range("L39").Value = 18000
$ws.Range("M39").Value = 18000
$ws.Range("P39").Value = 1000

$ws.Range("D40").Value2 = 44301
$ws.Range("I40").Value = "Especial"
$ws.Range("J40").Value = 30
$ws.Range("K40").Value = 18000
$ws.Range("L40").Value = 18000
$ws.Range("M40").Value = 18000
$ws.Range("P40").Value = 1000

$ws.Range("D41").Value2 = 44270
$ws.Range("J41").Value = 15
$ws.Range("K41").Value = 16000
$ws.Range("L41").Value = 16000
$ws.Range("M41").Value = 16000
$ws.Range("P41").Value = 889

$ws.Range("D42").Value2 = 44295
$ws.Range("J42").Value = 80
$ws.Range("K42").Value = 16000
$ws.Range("L42").Value = 16000
$ws.Range("M42").Value = 16000
$ws.Range("P42").Value = 889

$ws.Range("D43").Value2 = 44302
$ws.Range("J43").Value = 70
